$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

$ws1.Range("B1").Value = "2460-RBI-EPP-DB-SAR-REC-NOCOM-RNI-CTPD-SAR-MD-TR-1-DATE-VAR-INST-1st"
$ws1.Range("B2").Value = "246e"

$ws2.Range("B1").Value = "2460-RBI-EPP-DB-SAR-REC-NOCOM-RNI-CTPD-SAR-MD-TR-1-DATE-VAR-INST-1st"

$ws1.Activate() | Out-Null
$ws1.Range("B1").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("B1").Select() | Out-Null
